# Append one new data row (row 72) to each of the four sensor-data sheets,
# mirroring the existing row layout (A: timestamp text, B-E: hex text,
# F/H/I: numbers, G: big-integer ID kept as text).

$wb = $excel.ActiveWorkbook

$rowsToAdd = @(
    @{
        Sheet = "ROW35-FE-LIFTER"
        A = "2025-03-07 07:42:06"
        B = "0x01,0x90 "
        C = "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x14,0x41,0x0c,"
        D = "0x01,0x90,"
        E = "0x d"
        F = 400
        G = "568631262647113770877196"
        H = 400
        I = 13
    },
    @{
        Sheet = "ROW35-MID-LIFTER"
        A = "2025-03-07 07:29:35"
        B = "0x01,0x90 "
        C = "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x15,0x41,0x0c,"
        D = "0x01,0x90,"
        E = "0x e"
        F = 400
        G = "568631262647113770942732"
        H = 400
        I = 14
    },
    @{
        Sheet = "ROW02-FE-LIFTER"
        A = "2025-03-07 07:51:45"
        B = "0x01,0x90 "
        C = "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x06,0x41,0x0c,"
        D = "0x01,0x90,"
        E = "0xff"
        F = 400
        G = "568631262647113769959692"
        H = 400
        I = 255
    },
    @{
        Sheet = "ROW02-MID-LIFTER"
        A = "2025-03-07 07:41:15"
        B = "0x01,0x90 "
        C = "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x0b,0x40,0x0c,"
        D = "0x01,0x90,"
        E = "0x 3"
        F = 400
        G = "568631262647113769959692"
        H = 400
        I = 3
    }
)

foreach ($entry in $rowsToAdd) {
    $ws = $wb.Worksheets.Item($entry.Sheet)
    $newRow = $ws.UsedRange.Rows.Count + 1

    $ws.Cells.Item($newRow, 1).Value = $entry.A
    $ws.Cells.Item($newRow, 2).Value = $entry.B
    $ws.Cells.Item($newRow, 3).Value = $entry.C
    $ws.Cells.Item($newRow, 4).Value = $entry.D
    $ws.Cells.Item($newRow, 5).Value = $entry.E
    $ws.Cells.Item($newRow, 6).Value = $entry.F
    # Column G holds a 24-digit numeric identifier; enter it with a leading
    # apostrophe so Excel stores it as exact text instead of rounding it to
    # a floating-point double.
    $ws.Cells.Item($newRow, 7).Value = "'" + $entry.G
    $ws.Cells.Item($newRow, 8).Value = $entry.H
    $ws.Cells.Item($newRow, 9).Value = $entry.I
}
